# Insert a new weekly record as row 105, shifting the existing rows
# 105-127 down to 106-128 (matches the commit "Fruta / hortaliza, semanal").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 105 (and everything below it) down by one row.
$ws.Rows.Item(105).Insert()

# Populate the newly freed row 105 with the new weekly data point.
$ws.Cells.Item(105, 1).Value = 6
$ws.Cells.Item(105, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(105, 3).Value = "Metropolitana"
$ws.Cells.Item(105, 4).Value = 44951
$ws.Cells.Item(105, 5).Value = 13
$ws.Cells.Item(105, 6).Value = "Fruta"
$ws.Cells.Item(105, 7).Value = 100101
$ws.Cells.Item(105, 8).Value = "Berries"
$ws.Cells.Item(105, 9).Value = 100101008
$ws.Cells.Item(105, 10).Value = "Mora"
$ws.Cells.Item(105, 11).Value = "Sin especificar"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 200
$ws.Cells.Item(105, 14).Value = 4000
$ws.Cells.Item(105, 15).Value = 4000
$ws.Cells.Item(105, 16).Value = 4000
$ws.Cells.Item(105, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(105, 18).Value = "Región del Maule"
$ws.Cells.Item(105, 19).Value = 2000
$ws.Cells.Item(105, 20).Value = 2

# Make sure the D column keeps the date/time display format used by the
# rest of the column.
$ws.Cells.Item(105, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
